$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.825.44"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "3.458.04"
$ws.Range("E3").Value = "  +1.18%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'575.56"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'160.78"
$ws.Range("E6").Value = "  +2.36%  "

$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +13.29%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "3.463.72"
$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").Value = "'7.24"
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("E11").Value = "  +2.05%  "

$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  +3.89%  "

$ws.Range("D13").Value = "4.051.19"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "'28.31"
$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").Value = "64.890.49"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").Value = "3.490.79"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "'6.48"
$ws.Range("E19").Value = "  +3.13%  "

$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +2.20%  "

$ws.Range("D21").Value = "'381.33"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").Value = "'8.11"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").Value = "'0.552"
$ws.Range("E23").Value = "  +3.70%  "

$ws.Range("D24").Value = "'72.93"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").Value = "'10.08"
$ws.Range("E27").Value = "  +6.79%  "

$ws.Range("D28").Value = "'0.178"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").Value = "  +9.72%  "

$ws.Range("D31").Value = "'6.22"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("E32").Value = "  +1.51%  "

$ws.Range("D33").Value = "'23.56"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").Value = "'7.28"
$ws.Range("E34").Value = "  +6.64%  "

$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "  +11.59%  "

$ws.Range("D36").Value = "'160.95"
$ws.Range("E36").Value = "  +0.92%  "

$ws.Range("E37").Value = "  +4.75%  "

$ws.Range("D38").Value = "'0.0778"
$ws.Range("E38").Value = "  +2.84%  "

$ws.Range("D39").Value = "2.939.52"
$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("D40").Value = "'26.66"
$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").Value = "'6.73"
$ws.Range("E41").Value = "  +4.02%  "

$ws.Range("D42").Value = "'4.64"
$ws.Range("E42").Value = "  +7.18%  "

$ws.Range("D43").Value = "'0.0321"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("D44").Value = "'42.80"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("D45").Value = "'26.20"
$ws.Range("E45").Value = "  +12.34%  "

$ws.Range("D46").Value = "'0.777"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("E47").Value = "  +2.55%  "

$ws.Range("D48").Value = "'320.38"
$ws.Range("E48").Value = "  +9.98%  "

$ws.Range("E49").Value = "  +7.89%  "

$ws.Range("D50").Value = "'0.880"
$ws.Range("E50").Value = "  +5.19%  "

$ws.Range("D51").Value = "'2.18"
$ws.Range("E51").Value = "  -0.93%  "
